# Generate Report for Archive
#
# 1) Replace the "Ready for handoff" status text with "In Translation"
#    wherever it appears (Overview!E:F and the Status column, "C", on
#    each per-language sheet).
# 2) The Status column on each sheet was sized to fit its text; with the
#    new, shorter status string the columns narrow to match (mirrors
#    Excel re-fitting the "Status" column width after the handback
#    report is regenerated).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # NOTE: compare with the string literal on the LEFT of -eq. PowerShell's
        # -eq coerces the right-hand side to the left operand's type, so a
        # boolean-valued cell (text "True"/"False") compared as
        # ($cell.Value2 -eq "Ready for handoff") would coerce the literal to
        # [bool] (any non-empty string -> $true) and false-match every
        # "True" cell. Literal-first keeps this a plain string comparison.
        if ("Ready for handoff" -eq $cell.Value2) {
            $cell.Value = "In Translation"
        }
    }
}

# Narrow the "Status" columns to fit "In Translation" (was sized for the
# longer "Ready for handoff").
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.47   # column E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = 12.47   # column F (de-de status)

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.47        # column C (Status)

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.47        # column C (Status)
